$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.455.08"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.375.91"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.93%  "
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.981"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.93%  "
$ws.Range("D15").Value = "2.736.34"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").Value = "2.353.77"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").Value = "45.374.39"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.61%  "
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0955"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.24%  "
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.879.34"
$ws.Range("E45").Value = "  +12.81%  "
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.98%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.51%  "
